# "Criação do Manual do Usuário" — fills in the checklist answers on the
# "Ver-Transição1" sheet (IAP was previously #DIV/0! because the D column
# had no answers yet) and makes that sheet the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ver-Transição1")

# Checklist answers (D column). Values come from the sheet's data-validation
# list: "Sim", "Não", "NA".
$answers = @{
    "D6"  = "Sim"
    "D8"  = "Sim"
    "D10" = "Sim"
    "D11" = "Sim"
    "D12" = "NA"
    "D14" = "Sim"
    "D16" = "Sim"
    "D18" = "Sim"
    "D20" = "Não"
    "D21" = "Sim"
    "D22" = "Não"
    "D23" = "Não"
    "D24" = "Não"
    "D26" = "Sim"
    "D27" = "Não"
    "D28" = "Sim"
    "D30" = "Sim"
    "D31" = "Sim"
    "D33" = "Sim"
    "D34" = "Não"
    "D36" = "Sim"
    "D37" = "Sim"
    "D38" = "Sim"
    "D39" = "Sim"
    "D40" = "Sim"
    "D41" = "Sim"
    "D42" = "Não"
    "D44" = "Sim"
    "D45" = "Sim"
    "D46" = "Sim"
    "D48" = "Sim"
    "D49" = "Sim"
    "D50" = "Sim"
}

foreach ($addr in $answers.Keys) {
    $ws.Range($addr).Value = $answers[$addr]
}

# Recalculate so the IAP formula (F2) and every cell/chart depending on it
# (Indicadores!B6, the bar chart series) pick up the new answers instead of
# the old #DIV/0!.
$excel.CalculateFull()

# Make "Ver-Transição1" the active/selected sheet (tab moves from
# "Ver-Construção1" to "Ver-Transição1"), with D6 as the active cell,
# scrolled so row 31 is at the top — matching the saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("D6").Select()
